$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/tab to reflect the new extraction run (20240524-092026 -> 20240527-094031)
$ws.Name = "IClientBalance-20240527-094031-"

# Update the reference date column (G) for every data row from 45436 (2024-05-24) to 45439 (2024-05-27)
For ($r = 2; $r -le 257; $r++) {
    $ws.Cells.Item($r, 7).Value = 45439
}

# Correct the projected-value typos in rows 109 and 117 (D and H columns)
$ws.Cells.Item(109, 4).Value = 221.02
$ws.Cells.Item(109, 8).Value = 221.02

$ws.Cells.Item(117, 4).Value = 358.59
$ws.Cells.Item(117, 8).Value = 358.59
